$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 283
$ws.Range("C2").Value = 2016
$ws.Range("E2").Value = 82

$ws.Range("B3").Value = 462
$ws.Range("C3").Value = 519
$ws.Range("E3").Value = 3

$ws.Range("B4").Value = 241
$ws.Range("C4").Value = 289

$ws.Range("C5").Value = 2

$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 166
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 24
